$wb = $excel.ActiveWorkbook

# Update the "展览" (Exhibition) sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 149
$wsExhibit.Range("F4").Value = 101

# Update the "全部类型" (All types) sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 149
$wsAll.Range("F4").Value = 101
